$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows ("line7", "line8") right after the existing "line6"
# row (currently row 8) and push the "extr1".."extr8" rows down by two.
# ---------------------------------------------------------------------------
$ws.Rows("8:9").Insert()

# The inserted rows copy row 9's (old extr1, now pushed to row 10) format,
# so re-apply column A's formatting from a known-good sibling row (row 7,
# "line6") onto the two freshly inserted cells.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# New "line7" row (continues the branch index sequence: 6)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

# New "line8" row (branch index 7)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# ---------------------------------------------------------------------------
# The "extr1".."extr8" rows, now living at rows 10..17, keep their names but
# the branch-index column (A) continues the sequence (8..15) and a couple of
# rows get updated values.
# ---------------------------------------------------------------------------

$ws.Range("A10").Value = 8    # extr1
$ws.Range("A11").Value = 9    # extr2
$ws.Range("A12").Value = 10   # extr3
$ws.Range("A13").Value = 11   # extr4
$ws.Range("A14").Value = 12   # extr5
$ws.Range("A15").Value = 13   # extr6
$ws.Range("A16").Value = 14   # extr7
$ws.Range("A17").Value = 15   # extr8

# extr1 -> row 10 : in_service flips to TRUE
$ws.Range("E10").Value = $true

# extr2 -> row 11 : in_service flips to TRUE
$ws.Range("E11").Value = $true

# extr3 (row 12), extr4 (row 13), extr5 (row 14), extr6 (row 15),
# extr7 (row 16) and extr8 (row 17) keep their original C/D/E values -
# they simply shifted down by two rows thanks to the Insert above.

$wb.Save()
